$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation for Jengibre @ Vega Modelo de Temuco was
# recorded. It becomes the new first data row (row 28); every existing
# data row (28-107) shifts down by one (to 29-108).
$ws.Rows(28).Insert()

$ws.Range("A28").Value = 10
$ws.Range("B28").Value = "Vega Modelo de Temuco"
$ws.Range("C28").Value = "La Araucanía"
$ws.Range("D28").Value = 44453
$ws.Range("E28").Value = 9
$ws.Range("F28").Value = 100114007
$ws.Range("G28").Value = "Jengibre"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 40
$ws.Range("K28").Value = 20000
$ws.Range("L28").Value = 25000
$ws.Range("M28").Value = 22500
$ws.Range("N28").Value = "$/caja 13 kilos"
$ws.Range("O28").Value = "Perú"
$ws.Range("P28").Value = 1731
$ws.Range("Q28").Value = 13
$ws.Range("R28").Value = "Hortaliza"
